$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new blank rows before the old "Total Hours" row (old row 75),
# pushing it down to row 81 and making room for new log rows 74-78.
$ws.Rows("75:80").Insert()

# Copy formatting (styles, number formats, alignment) from the last
# existing data row (73) down onto the 5 new rows (74-78).
$ws.Range("A73:G73").Copy()
$ws.Range("A74:G78").PasteSpecial(-4122)

# --- Row 74 ---
$ws.Range("A74").Value = 73
$ws.Range("B74").Value = 44766
$ws.Range("C74").Value = 0.51041666666666663
$ws.Range("D74").Value = 0.64583333333333337
$ws.Range("E74").Formula = "=D74-C74"
$ws.Range("F74").Value = "Doc"
$ws.Range("G74").Value = "1. PPT and block diagram for PSPNet, FCN and Unet "

# --- Row 75 ---
$ws.Range("A75").Value = 74
$ws.Range("B75").Value = 44767
$ws.Range("C75").Value = 0.33333333333333331
$ws.Range("D75").Value = 0.38541666666666669
$ws.Range("E75").Formula = "=D75-C75"
$ws.Range("F75").Value = "Doc"
$ws.Range("G75").Value = "1. PSPNet block diagram and PPT completed"

# --- Row 76 ---
$ws.Range("A76").Value = 75
$ws.Range("B76").Value = 44767
$ws.Range("C76").Value = 0.41666666666666669
$ws.Range("D76").Value = 0.59375
$ws.Range("E76").Formula = "=D76-C76"
$ws.Range("F76").Value = "Doc"
$ws.Range("G76").Value = "1. deeplabv3+, mix FFN block diagrams added"

# --- Row 77 ---
$ws.Range("A77").Value = 76
$ws.Range("B77").Value = 44767
$ws.Range("C77").Value = 0.63541666666666663
$ws.Range("D77").Value = 0.75
$ws.Range("E77").Formula = "=D77-C77"
$ws.Range("F77").Value = "Doc"
$ws.Range("G77").Value = "1. Dot product attention and multi head self attention block diagrams added"
$ws.Rows("77:77").RowHeight = 30

# --- Row 78 ---
$ws.Range("A78").Value = 77
$ws.Range("B78").Value = 44767
$ws.Range("C78").Value = 0.84375
$ws.Range("D78").Value = 0.875
$ws.Range("E78").Formula = "=D78-C78"
$ws.Range("F78").Value = "Doc"
$ws.Range("G78").Value = "1. References for Libraries, research papers, code documented"

# --- Update the Total Hours formula which is now on row 81 ---
$ws.Range("E81").Formula = "=SUM(E2:E80)"

# --- Update view selection to match the end of the edit session ---
$ws.Range("G79").Select()
